$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bumped by one day
# (2026-02-28 -> 2026-03-01, i.e. serial 46081 -> 46082) for every data row
# (rows 2 through 276).
$ws.Range("C2:C276").Value = 46082
